$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New reconciliation rows (86-100) appended below the existing data.
$newRows = @(
    @(237671823369, "MFS ENTREE COLLEGE MALANGUE", "Rte_1", "Cite Sic Stade Marion Ocm", 189785.8, 157448, -32337.79999999999, 0.8296089591528977, "Cite Sic"),
    @(237672128028, "CAROLINE WAKO DJAMNOU", "Rte_1", "Cite Sic Stade Marion Ocm", 25030, 65097, 40067, 2.600759089093089, "Cite Sic"),
    @(237672277367, "TOP MOBIL KM5 LTDLA_POLAS_BTQ_KM5", "Rte_0", "Cite Sic Stade Marion Ocm", 26250, 2996, -23254, 0.1141333333333333, "Cite Sic"),
    @(237674853971, "NJOSSEU TCHOUNZOU TOP MOBILE", "Rte_1", "Cite Sic Stade Marion Ocm", 182169.7, 4306, -177863.7, 0.02363730082445105, "Cite Sic"),
    @(237674884705, "BAH AMADOU MOUNTAGHA ETS MOBILE FINANCIAL SERVICES MFS", "Rte_0", "Cite Sic Stade Marion Ocm", 5000, 9333, 4333, 1.8666, "Cite Sic"),
    @(237675779272, "RODES NGWEM KEMAYOU", "Rte_1", "Cite Sic Stade Marion Ocm", 184162.0370370371, 30585, -153577.0370370371, 0.1660765730661907, "Cite Sic"),
    @(237677304210, "CARINE OROCK", "Rte_2", "Cite Sic Stade Marion Ocm", 61193.33333333334, 235659, 174465.6666666667, 3.851056759995642, "Cite Sic"),
    @(237678267353, "LA NEGRESSE SARL EMBOLA BELTUS MBU", "Rte_0", "Cite Sic Stade Marion Ocm", 16520, 0, -16520, 0, "Cite Sic"),
    @(237678370615, "ESSEN ONGOLONG BERTHE HORTENSE ETS MOBILE FINANCIAL SERVICES MFS", "Rte_1", "Cite Sic Stade Marion Ocm", 100229.8, 510310, 410080.2, 5.09139996288529, "Cite Sic"),
    @(237678836319, "KAMDOM DOMINIQUE STEPHANIE ETS MOBILE FINANCIAL SERVICES MFS", "Rte_1", "Cite Sic Stade Marion Ocm", 25619.09090909091, 137271, 111651.9090909091, 5.358152656044853, "Cite Sic"),
    @(237678922502, "NWOAGA TCHAMDJOU EPSE KAMSEU EMILINE ETS LE CONTENT", "Rte_0", "Cite Sic Stade Marion Ocm", 56220, 730908, 674688, 13.0008537886873, "Cite Sic"),
    @(237679884264, "GABRIEL MONKAM TCHOUPE", "Rte_1", "Cite Sic Stade Marion Ocm", 59610, 365895, 306285, 6.138147961751384, "Cite Sic"),
    @(237681019523, "ETS MOULAY RIPERT AND COMPANY", "Rte_1", "Cite Sic Stade Marion Ocm", 57805.71428571429, 245762, 187956.2857142857, 4.251517398181099, "Cite Sic"),
    @(237681125655, "EMENGUE PICHOU ROMEO KAMILAH CONNECTION GROUP", "Rte_1", "Cite Sic Stade Marion Ocm", 47813.75, 33916, -13897.75, 0.7093357036417348, "Cite Sic"),
    @(237681240793, "MBANE EMILIE FRANCOISE ETS MOBILE FINANCIAL SERVICES MFS", "Rte_0", "Cite Sic Stade Marion Ocm", 7438.090909090909, 2067, -5371.090909090909, 0.2778938877277894, "Cite Sic"),
)

$startRow = 86
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowData[$c]
    }
}
